$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column C
$ws.Range("C1").Value = "Abbrev"
$ws.Range("C1").Font.Bold = $true

# Region / country abbreviation codes, in the same row order as column B (rows 2-50)
$codes = @(
    "AT","BE","BG","CY","CZ","DE","DK","EE","ES","FI",
    "FR","GR","HR","HU","IE","IT","LT","LU","LV","MT",
    "NL","PL","PT","RO","SE","SI","SK","GB","US","JP",
    "CN","CA","KR","BR","IN","MX","RU","AU","CH","TR",
    "TW","NO","ID","ZA","WA","WL","WE","WF","WM"
)

for ($i = 0; $i -lt $codes.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $codes[$i]
}

# Match print setup changes captured in the source workbook
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Match the selection state recorded in the source workbook
$ws.Range("E5").Select() | Out-Null
